$wb = $excel.ActiveWorkbook
$styles = $wb.Styles
$styles | Get-Member -MemberType Method | Out-String | Write-Host
